$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.188.98"
$ws.Range("E2").Value = "  +0.67%  "
$ws.Range("D3").Value = "3.502.48"
$ws.Range("E3").Value = "  -0.05%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "604.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.89"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.11%  "
$ws.Range("E7").Value = "  -0.67%  "
$ws.Range("D8").Value = "3.502.75"
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("E10").Value = "  -2.37%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.19"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +7.29%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.587"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.99%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "46.04"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.15%  "
$ws.Range("E14").Value = "  -1.09%  "
$ws.Range("D15").Value = "4.078.70"
$ws.Range("E15").Value = "  +0.11%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "614.10"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.07%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "8.27"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.35%  "
$ws.Range("D18").Value = "3.507.92"
$ws.Range("E18").Value = "  +0.16%  "
$ws.Range("D19").Value = "70.254.01"
$ws.Range("E19").Value = "  +0.79%  "
$ws.Range("E20").Value = "  +1.03%  "
$ws.Range("E21").Value = "  +1.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.878"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.42%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.05"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -8.49%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "99.26"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.31%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "15.56"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.43%  "
$ws.Range("E26").Value = "  -3.21%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.56"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.19%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "34.01"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.86%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.01"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.44%  "
$ws.Range("B31").Value = "Bittensor"
$ws.Range("C31").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "659.93"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +17.43%  "
$ws.Range("B32").Value = "Stacks"
$ws.Range("C32").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.97"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.69%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "8.03"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.19%  "
$ws.Range("E34").Value = "  -4.41%  "
$ws.Range("E35").Value = "  -2.15%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0995"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.40%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.71"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.26%  "
$ws.Range("E38").Value = "  +0.49%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0475"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.30%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "56.63"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.61%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.144"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.53%  "
$ws.Range("D43").Value = "3.350.79"
$ws.Range("E43").Value = "  +0.85%  "
$ws.Range("E44").Value = "  +3.87%  "
$ws.Range("E45").Value = "  -4.94%  "
$ws.Range("E46").Value = "  -2.46%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "31.86"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.13%  "
$ws.Range("E48").Value = "  -2.42%  "
$ws.Range("E49").Value = "  +0.75%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "132.79"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.60%  "
$ws.Range("E51").Value = "  -0.02%  "
